$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.511.86'
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").Value = '1.903.88'
$ws.Range("E3").Value = '  +1.50%  '

$ws.Range("D4").Value = '''0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '''238.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("D6").Value = '''0.9989'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").Value = '''0.4910'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.73%  '

$ws.Range("D8").Value = '''0.2935'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.27%  '

$ws.Range("D9").Value = '''0.06709'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("D10").Value = '1.903.41'
$ws.Range("E10").Value = '  +1.56%  '

$ws.Range("D11").Value = '''17.03'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.01%  '

$ws.Range("D12").Value = '''0.07340'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.41%  '

$ws.Range("D13").Value = '''5.189'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.91%  '

$ws.Range("D14").Value = '''88.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '

$ws.Range("D15").Value = '''0.6710'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.28%  '

$ws.Range("D16").Value = '30.479.10'
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").Value = '''13.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.23%  '

$ws.Range("D18").Value = '''0.000007902'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.71%  '

$ws.Range("D19").Value = '''0.9995'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.10%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''5.486'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +16.60%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.137.61'
$ws.Range("E21").Value = '  +0.95%  '

$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").Value = '''195.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.30%  '

$ws.Range("D24").Value = '''6.150'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.87%  '

$ws.Range("D25").Value = '''9.523'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.18%  '

$ws.Range("D26").Value = '''163.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.15%  '

$ws.Range("D27").Value = '''18.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("D28").Value = '''1.953'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.24%  '

$ws.Range("D29").Value = '''1.470'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.62%  '

$ws.Range("D30").Value = '''4.356'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.53%  '

$ws.Range("D31").Value = '''0.09193'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.81%  '

$ws.Range("D32").Value = '''4.089'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.50%  '

$ws.Range("D33").Value = '''0.05175'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.31%  '

$ws.Range("D34").Value = '''0.7463'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.58%  '

$ws.Range("D35").Value = '''1.109'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.96%  '

$ws.Range("D36").Value = '''2.718'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.04%  '

$ws.Range("D37").Value = '''0.01819'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("E38").Value = '  +1.12%  '

$ws.Range("D39").Value = '''0.9261'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.19%  '

$ws.Range("D40").Value = '''2.065'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.94%  '

$ws.Range("D41").Value = '''0.4402'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.39%  '

$ws.Range("D42").Value = '''106.91'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.22%  '

$ws.Range("D43").Value = '''5.929'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.85%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '''69.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +22.13%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''0.9954'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("D46").Value = '''0.1376'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.77%  '

$ws.Range("D47").Value = '''7.620'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.36%  '

$ws.Range("D48").Value = '''9.053'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.13%  '

$ws.Range("D49").Value = '''35.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.15%  '

$ws.Range("D50").Value = '''0.05847'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.55%  '

$ws.Range("D51").Value = '''0.3934'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.19%  '
